# Auto-generated: apply updated market price data to each profession sheet.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 102.15385
$ws.Range("J2").Value = 157.33333
$ws.Range("L2").Value = 157.33333
$ws.Range("N2").Value = -383.33333
$ws.Range("H80").Value = 527621
$ws.Range("I80").Value = 867.44446
$ws.Range("J80").Value = 1001699.2
$ws.Range("K80").Value = 2602.33338
$ws.Range("L80").Value = 3005097.6
$ws.Range("M80").Value = -1604.33338
$ws.Range("N80").Value = -3007093.6
$ws.Range("H83").Value = 527621
$ws.Range("I83").Value = 867.44446
$ws.Range("J83").Value = 1001699.2
$ws.Range("K83").Value = 7807.00014
$ws.Range("L83").Value = 9015292.799999999
$ws.Range("M83").Value = -2815.00014
$ws.Range("N83").Value = -9025276.799999999
$ws.Range("H88").Value = 2319.6
$ws.Range("I88").Value = 2280.4
$ws.Range("K88").Value = 2280.4
$ws.Range("M88").Value = -1874.4
$ws.Range("H91").Value = 2319.6
$ws.Range("I91").Value = 2280.4
$ws.Range("K91").Value = 2280.4
$ws.Range("M91").Value = -876.4000000000001
$ws.Range("H98").Value = 2459.5386
$ws.Range("J98").Value = 3665.6667
$ws.Range("L98").Value = 3665.6667
$ws.Range("N98").Value = -6661.6667
$ws.Range("H111").Value = 1189.1
$ws.Range("I111").Value = 1315.5883
$ws.Range("J111").Value = 472.33334
$ws.Range("K111").Value = 3946.7649
$ws.Range("L111").Value = 1417.00002
$ws.Range("M111").Value = -879.7648999999997
$ws.Range("N111").Value = -7551.000019999999
$ws.Range("H122").Value = 2459.5386
$ws.Range("J122").Value = 3665.6667
$ws.Range("L122").Value = 10997.0001
$ws.Range("N122").Value = -15897.0001
$ws.Range("H127").Value = 2498.3333
$ws.Range("I127").Value = 2497.5
$ws.Range("K127").Value = 7492.5
$ws.Range("M127").Value = -2532.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2855.8286
$ws.Range("I2").Value = 2144.6086
$ws.Range("K2").Value = 2144.6086
$ws.Range("M2").Value = -2031.6086
$ws.Range("H32").Value = 4088.3416
$ws.Range("I32").Value = 3635
$ws.Range("K32").Value = 3635
$ws.Range("M32").Value = -3348
$ws.Range("H61").Value = 11076.655
$ws.Range("I61").Value = 9815.066000000001
$ws.Range("J61").Value = 12428.357
$ws.Range("K61").Value = 9815.066000000001
$ws.Range("L61").Value = 12428.357
$ws.Range("M61").Value = -9603.066000000001
$ws.Range("N61").Value = -12852.357
$ws.Range("H105").Value = 52184.5
$ws.Range("J105").Value = 52184.5
$ws.Range("L105").Value = 52184.5
$ws.Range("N105").Value = -59172.5
$ws.Range("H116").Value = 2855.8286
$ws.Range("I116").Value = 2144.6086
$ws.Range("K116").Value = 2144.6086
$ws.Range("M116").Value = 149.3914
$ws.Range("H136").Value = 11076.655
$ws.Range("I136").Value = 9815.066000000001
$ws.Range("J136").Value = 12428.357
$ws.Range("K136").Value = 29445.198
$ws.Range("L136").Value = 37285.071
$ws.Range("M136").Value = -26895.198
$ws.Range("N136").Value = -42385.071

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2855.8286
$ws.Range("I3").Value = 2144.6086
$ws.Range("K3").Value = 2144.6086
$ws.Range("M3").Value = -2030.6086
$ws.Range("H94").Value = 1963.875
$ws.Range("J94").Value = 2628.8333
$ws.Range("L94").Value = 2628.8333
$ws.Range("N94").Value = -3530.8333
$ws.Range("H99").Value = 2021.6487
$ws.Range("I99").Value = 2471.6538
$ws.Range("K99").Value = 2471.6538
$ws.Range("M99").Value = -973.6538
$ws.Range("H103").Value = 33100.25
$ws.Range("J103").Value = 33100.25
$ws.Range("L103").Value = 33100.25
$ws.Range("N103").Value = -35444.25

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2520.25
$ws.Range("I16").Value = 2737.647
$ws.Range("J16").Value = 1992.2858
$ws.Range("K16").Value = 2737.647
$ws.Range("L16").Value = 1992.2858
$ws.Range("M16").Value = -2450.647
$ws.Range("N16").Value = -2566.2858
$ws.Range("H86").Value = 4166.9443
$ws.Range("I86").Value = 2864.2856
$ws.Range("K86").Value = 2864.2856
$ws.Range("M86").Value = -1741.2856
$ws.Range("H89").Value = 4166.9443
$ws.Range("I89").Value = 2864.2856
$ws.Range("K89").Value = 14321.428
$ws.Range("M89").Value = -8705.428
$ws.Range("H99").Value = 10423.914
$ws.Range("I99").Value = 5724.8887
$ws.Range("J99").Value = 12050.5
$ws.Range("K99").Value = 5724.8887
$ws.Range("L99").Value = 12050.5
$ws.Range("M99").Value = -4226.8887
$ws.Range("N99").Value = -15046.5
$ws.Range("H100").Value = 780
$ws.Range("J100").Value = 780
$ws.Range("L100").Value = 780
$ws.Range("N100").Value = -2944
$ws.Range("H113").Value = 2520.25
$ws.Range("I113").Value = 2737.647
$ws.Range("J113").Value = 1992.2858
$ws.Range("K113").Value = 2737.647
$ws.Range("L113").Value = 1992.2858
$ws.Range("M113").Value = -567.6469999999999
$ws.Range("N113").Value = -6332.2858
$ws.Range("H121").Value = 39993.25
$ws.Range("J121").Value = 39993.25
$ws.Range("L121").Value = 39993.25
$ws.Range("N121").Value = -42613.25
$ws.Range("H126").Value = 10423.914
$ws.Range("I126").Value = 5724.8887
$ws.Range("J126").Value = 12050.5
$ws.Range("K126").Value = 17174.6661
$ws.Range("L126").Value = 36151.5
$ws.Range("M126").Value = -14704.6661
$ws.Range("N126").Value = -41091.5
$ws.Range("H132").Value = 21334.334
$ws.Range("I132").Value = 12661.648
$ws.Range("K132").Value = 37984.944
$ws.Range("M132").Value = -35454.944

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3015.7144
$ws.Range("I131").Value = 970
$ws.Range("J131").Value = 3460.4348
$ws.Range("K131").Value = 2910
$ws.Range("L131").Value = 10381.3044
$ws.Range("M131").Value = 2130
$ws.Range("N131").Value = -20461.3044
$ws.Range("H133").Value = 5659.8
$ws.Range("J133").Value = 6400
$ws.Range("L133").Value = 19200
$ws.Range("N133").Value = -29320
$ws.Range("H140").Value = 2049.375
$ws.Range("I140").Value = 2049.375
$ws.Range("K140").Value = 6148.125
$ws.Range("M140").Value = -968.125

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 74.55556
$ws.Range("I2").Value = 42.666668
$ws.Range("J2").Value = 90.5
$ws.Range("K2").Value = 42.666668
$ws.Range("L2").Value = 90.5
$ws.Range("M2").Value = 70.333332
$ws.Range("N2").Value = -316.5
$ws.Range("H70").Value = 9985.166999999999
$ws.Range("J70").Value = 9995.666999999999
$ws.Range("L70").Value = 9995.666999999999
$ws.Range("N70").Value = -10535.667
$ws.Range("H73").Value = 9985.166999999999
$ws.Range("J73").Value = 9995.666999999999
$ws.Range("L73").Value = 9995.666999999999
$ws.Range("N73").Value = -11867.667
$ws.Range("H92").Value = 31530.166
$ws.Range("J92").Value = 30941.555
$ws.Range("L92").Value = 30941.555
$ws.Range("N92").Value = -34685.555
$ws.Range("H113").Value = 79104.34
$ws.Range("J113").Value = 2309.8333
$ws.Range("L113").Value = 2309.8333
$ws.Range("N113").Value = -6649.8333
$ws.Range("H132").Value = 12487.049
$ws.Range("I132").Value = 10500.238
$ws.Range("J132").Value = 16659.35
$ws.Range("K132").Value = 31500.714
$ws.Range("L132").Value = 49978.05
$ws.Range("M132").Value = -28970.714
$ws.Range("N132").Value = -55038.05

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8071.909
$ws.Range("I16").Value = 8728.200000000001
$ws.Range("J16").Value = 6665.5713
$ws.Range("K16").Value = 8728.200000000001
$ws.Range("L16").Value = 6665.5713
$ws.Range("M16").Value = -8558.200000000001
$ws.Range("N16").Value = -7005.5713
$ws.Range("H82").Value = 1445.4
$ws.Range("I82").Value = 1329.3334
$ws.Range("J82").Value = 1561.4667
$ws.Range("K82").Value = 1329.3334
$ws.Range("L82").Value = 1561.4667
$ws.Range("M82").Value = -968.3334
$ws.Range("N82").Value = -2283.4667
$ws.Range("H85").Value = 1445.4
$ws.Range("I85").Value = 1329.3334
$ws.Range("J85").Value = 1561.4667
$ws.Range("K85").Value = 1329.3334
$ws.Range("L85").Value = 1561.4667
$ws.Range("M85").Value = -81.33339999999998
$ws.Range("N85").Value = -4057.4667
$ws.Range("H122").Value = 3813.762
$ws.Range("I122").Value = 4240
$ws.Range("K122").Value = 12720
$ws.Range("M122").Value = -10270
$ws.Range("H132").Value = 5313.4
$ws.Range("I132").Value = 5019.696
$ws.Range("J132").Value = 5876.3335
$ws.Range("K132").Value = 15059.088
$ws.Range("L132").Value = 17629.0005
$ws.Range("M132").Value = -12529.088
$ws.Range("N132").Value = -22689.0005

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 26155.77
$ws.Range("J54").Value = 26155.77
$ws.Range("L54").Value = 26155.77
$ws.Range("N54").Value = -27195.77
$ws.Range("H98").Value = 22063.334
$ws.Range("J98").Value = 22063.334
$ws.Range("L98").Value = 22063.334
$ws.Range("N98").Value = -28053.334
$ws.Range("H103").Value = 17200.666
$ws.Range("J103").Value = 17200.666
$ws.Range("L103").Value = 17200.666
$ws.Range("N103").Value = -19544.666
$ws.Range("H107").Value = 11114179
$ws.Range("I107").Value = 1580.5
$ws.Range("J107").Value = 27783078
$ws.Range("K107").Value = 4741.5
$ws.Range("L107").Value = 83349234
$ws.Range("M107").Value = -2821.5
$ws.Range("N107").Value = -83353074
$ws.Range("H132").Value = 148958.6
$ws.Range("I132").Value = 241905.06
$ws.Range("J132").Value = 24061.75
$ws.Range("K132").Value = 725715.1799999999
$ws.Range("L132").Value = 72185.25
$ws.Range("M132").Value = -723185.1799999999
$ws.Range("N132").Value = -77245.25
$ws.Range("H136").Value = 2944776.2
$ws.Range("I136").Value = 4258962
$ws.Range("K136").Value = 12776886
$ws.Range("M136").Value = -12774336
